# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update case counts for a couple of countries (Kazajistan, Australia)
# - Re-sort three pairs of countries whose "Casos activos" (active cases)
#   changed enough to swap their ranking order (Belice/Nueva Caledonia,
#   Papua Nueva Guinea/Islas Virgenes Britanicas, San Bartolome/Bonaire)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 05:05"

# --- Plain value refreshes (no re-sort) --------------------------------
# Row 54: Kazajistan
$ws.Range("B54").Value = 9576
$ws.Range("C54").Value = 272
$ws.Range("E54").Value = 4771

# Row 66: Australia
$ws.Range("B66").Value = 7150
$ws.Range("C66").Value = 11
$ws.Range("D66").Value = 6579
$ws.Range("E66").Value = 468

# --- Re-sorted pairs ----------------------------------------------------
# Rows 200/201: Belice <-> Nueva Caledonia swap places
$ws.Range("A200").Value = "Nueva Caledonia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# Rows 213/214: Papua Nueva Guinea <-> Islas Virgenes Britanicas swap places
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0

# Rows 215/216: San Bartolome <-> Bonaire, San Eustaquio y Saba swap places
# (values are identical between the two, only the labels move)
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"
